# Insert a new "most recent quarter" column before column D on the SIM sheet,
# shifting the existing quarterly columns right, then refresh all the figures
# with the restated financial data (new quarter added, several trailing
# quarters restated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column; Excel shifts D:K -> E:L. The new column inherits the
# formatting of column C, so copy the (bounded, used-range only) formatting
# from column E -- which still carries the original column D formatting --
# into the new column D. Using a bounded range (not whole columns) keeps the
# sheet's used range/dimension intact.
$ws.Columns("D").Insert()
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Row => new values for columns D:L (9 quarterly columns incl. the new one)
$rowData = @{
    7 = @(43373,43281,43190,43100,43008,42916,42825,42735,42643)
    8 = @(509900,880400,411100,352100,352200,370300,409800,424900,354200)
    9 = @(435500,726400,327900,302100,302900,307700,328300,365400,288600)
    10 = @(74400,154000,83300,50000,49300,62600,81400,59500,65600)
    11 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    12 = @("NA","NA","NA","NA","NA","NA","NA","NA","NA")
    13 = @(0,0,0,0,0,0,0,0,0)
    14 = @("NA","NA","NA","NA","NA","NA","NA",0,"NA")
    15 = @(0,0,0,0,0,0,0,0,0)
    16 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    17 = @(449300,754000,343200,315500,314400,326500,348400,396200,301800)
    18 = @(60600,126400,68000,36600,37900,43800,61400,28600,52400)
    19 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    20 = @(20900,-23700,-15600,-1600,8200,-10900,-26800,34600,22100)
    21 = @(95800,129900,66200,54500,63600,51300,55000,91300,92200)
    22 = @(1300,3500,600,900,700,600,600,1100,500)
    23 = @(80100,99200,51800,34200,45400,32300,33900,62200,74100)
    24 = @(5500,7400,500,58100,1400,600,-3300,52500,-3100)
    25 = @(0,0,0,0,0,0,0,0,0)
    26 = @(74600,91700,51300,-23900,44000,31700,37200,9700,77100)
    27 = @(71900,92900,50400,-16000,38100,31700,35200,-73200,77600)
    28 = @(0,0,0,0,0,0,0,0,0)
    29 = @(0,0,0,0,0,0,0,0,0)
    30 = @(0,0,0,0,0,0,0,0,0)
    31 = @(0,0,0,0,0,0,0,0,0)
    32 = @(-20900,23700,15600,1600,-8200,10900,26800,-34600,-22100)
    33 = @(71900,92900,50400,-16000,38100,31700,35200,-73200,77600)
    34 = @(0,0,0,0,0,0,0,0,0)
    35 = @(71900,92900,50400,-16000,38100,31700,35200,-73200,77600)
    38 = @(43373,43281,43190,43100,43008,42916,42825,42735,42643)
    39 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    40 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    41 = @(343800,386400,350700,377200,343500,343400,351000,422900,406000)
    42 = @(10100,5900,5100,4400,7700,4900,3500,3900,6100)
    43 = @(635200,575400,557800,519400,482600,440500,450800,553200,339400)
    44 = @(521400,472600,409300,428400,402800,383300,336000,529300,331700)
    45 = @(37900,37000,28000,22000,25800,29100,36200,40700,17000)
    46 = @(1548400,1477300,1351000,1351400,1262500,1201100,1177400,1221700,1100300)
    47 = @(0,0,0,0,0,0,0,0,0)
    48 = @(832100,843800,798800,813900,773900,742700,749200,1538900,740500)
    49 = @(122000,124900,127100,127800,128000,127600,130200,270200,140000)
    50 = @(0,0,0,0,0,0,0,0,0)
    51 = @(0,0,0,0,0,0,0,0,0)
    52 = @(67000,62500,58300,62100,73700,73100,102200,84800,88000)
    53 = @(0,0,0,0,0,0,0,0,0)
    54 = @(2569500,2508500,2335200,2355300,2238200,2144600,2159100,2211100,2068700)
    55 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    56 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    57 = @(273600,216100,207000,284500,191700,178900,174500,150100,162400)
    58 = @(300,300,300,300,300,300,300,300,300)
    59 = @(232800,218800,160700,89500,155700,141400,152100,195600,149900)
    60 = @(506800,435200,368000,374300,347600,320600,326900,293100,312600)
    61 = @(0,0,0,0,0,0,0,0,0)
    62 = @(233300,252800,206500,216100,149300,136700,90900,181500,100300)
    63 = @(0,0,0,0,0,0,0,0,0)
    64 = @(0,0,0,0,0,0,0,0,0)
    65 = @(0,0,0,0,0,0,0,0,0)
    66 = @(736000,681000,569900,584800,499300,454500,378900,442000,321400)
    67 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    68 = @(0,0,0,0,0,0,0,0,0)
    69 = @(0,0,0,0,0,0,0,0,0)
    70 = @(0,0,0,0,0,0,0,0,0)
    71 = @(0,0,0,0,0,0,0,0,0)
    72 = @(1084100,1005000,1124900,1074500,1090500,1052400,1107800,1117300,1138100)
    73 = @(0,0,0,0,0,0,0,0,0)
    74 = @(0,0,0,0,0,0,0,0,0)
    75 = @(0,0,0,0,0,0,0,0,0)
    76 = @(1833500,1827500,1765300,1770500,1738800,1690100,1780200,1769000,1747300)
    77 = @(0,0,0,0,0,0,0,0,0)
    80 = @(43373,43281,43190,43100,43008,42916,42825,42735,42643)
    81 = @(71900,92900,50400,-16000,38100,31700,35200,-73200,77600)
    82 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    83 = @(14400,27200,13900,19400,17600,18400,20500,28100,17700)
    84 = @(0,0,0,0,0,0,0,0,0)
    85 = @(0,0,0,0,0,0,0,0,0)
    86 = @(0,0,0,0,0,0,0,0,0)
    87 = @(0,0,0,0,0,0,0,0,0)
    88 = @(0,0,0,0,0,0,0,0,0)
    89 = @(6500,126900,-6800,70700,42200,31400,-900,3600,137900)
    90 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    91 = @(-19000,-37500,-18500,-40000,-40300,-37900,-39000,-42900,-52600)
    92 = @(0,0,0,0,0,0,0,0,0)
    93 = @(0,0,0,0,0,0,0,0,0)
    94 = @(-15000,-117700,-18000,-21200,-40000,-37100,-41700,-42000,-46500)
    95 = @($null,$null,$null,$null,$null,$null,$null,$null,$null)
    96 = @(0,0,0,0,0,0,0,0,0)
    97 = @(0,0,0,0,0,0,0,0,0)
    98 = @(0,0,0,0,0,0,0,0,0)
    99 = @(0,0,0,0,0,0,0,0,0)
    100 = @(-16000,-4800,0,-4400,-2100,-3700,-9100,49800,-3600)
    101 = @(-18000,4800,-1600,-11500,0,-2700,-4600,2600,-46300)
    102 = @(-42600,9300,-26400,33600,100,-8700,-59800,14000,41600)
}

foreach ($rowNum in $rowData.Keys) {
    $vals = $rowData[$rowNum]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $cell = $ws.Cells.Item([int]$rowNum, 4 + $i)
        if ($null -eq $vals[$i]) {
            $cell.Value2 = $null
        } else {
            $cell.Value2 = $vals[$i]
        }
    }
}
